# Update "想去人数" (column F) values on sheets "展览" and "全部类型"
# to reflect newly generated counts (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 64
$ws1.Range("F3").Value = 11635
$ws1.Range("F4").Value = 212
$ws1.Range("F5").Value = 333
$ws1.Range("F7").Value = 11609
$ws1.Range("F8").Value = 480
$ws1.Range("F10").Value = 88
$ws1.Range("F11").Value = 1763
$ws1.Range("F12").Value = 5744
$ws1.Range("F14").Value = 3509
$ws1.Range("F16").Value = 16

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 64
$ws4.Range("F5").Value = 11635
$ws4.Range("F6").Value = 212
$ws4.Range("F7").Value = 333
$ws4.Range("F9").Value = 11609
$ws4.Range("F10").Value = 480
$ws4.Range("F12").Value = 88
$ws4.Range("F13").Value = 1763
$ws4.Range("F15").Value = 5744
$ws4.Range("F17").Value = 3509
$ws4.Range("F19").Value = 16
